$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Mid" simulation results (B2:K11), re-run after modifying the
# simulation constraints; values are now plain whole numbers instead of the
# previous long decimal expansions.
$rows = @(
  @(787,351,206,134,90,61,40,24,10,0),
  @(787,351,206,134,90,61,40,24,11,11),
  @(787,352,206,134,90,61,40,24,24,24),
  @(788,352,207,134,90,61,40,42,42,42),
  @(788,352,207,134,90,61,63,65,66,66),
  @(788,352,207,134,90,93,96,99,99,99),
  @(789,352,207,134,138,142,145,148,149,149),
  @(789,352,207,212,216,216,216,216,216,216),
  @(532,216,216,216,216,216,216,216,216,216),
  @(216,216,216,216,216,216,216,216,216,216)
)

$arr = New-Object 'object[,]' 10,10
for ($r = 0; $r -lt 10; $r++) {
  for ($c = 0; $c -lt 10; $c++) {
    $arr[$r, $c] = $rows[$r][$c]
  }
}

$ws.Range("B2:K11").Value = $arr

# Add the (currently empty) number-formatted cells that will hold the data
# for the new "adding a nodule to a plant that has/doesn't have AMF" figure,
# columns M:V, rows 2:11.
$ws.Range("M2:V11").NumberFormat = "0"

# Restore the selection left on the sheet when the author last saved it.
$ws.Range("A14:XFD25").Select()
